$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "'281474991205821-1748891368335"
$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Value = "Mobile Usage"
$ws.Range("C2").Value = "'2025-06-02T13:09:28.335"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").Value = "'281474991205821"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'148"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "'51834015"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "LUIS IBARRA"
$ws.Range("H2").Value = 20.514490729
$ws.Range("I2").Value = -103.36359083
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1748891365835/F7KsVXPM8Y-camera-video-segment-driver-1748891368335.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSORYA3UK6%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T151442Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjED8aCXVzLXdlc3QtMiJGMEQCIA0KiZttqleQAYFkoPrqGfaQ%2BEfVrMoOu8tcRRo3a5ubAiAR1SS6ywof5CKQTSha4d49IXbz7Aet18sLqGw7V0KBHCrdAwgYEAQaDDc4MTIwNDk0MjI0NCIMKvKB7UOIkRyDvF9OKroDViZbjkBBAv%2FiogeHYKxF%2BDalm8ruSrGzUNAdMd1bL85jwGiCmx%2BvYVaoGRUS6DRdywOkDnUGtHrAWIFtDmLdQgAelEr%2FMOD8N0Z3%2FcKX5oi3cIHVXVh87O5WXRp2Y5kQwrzy4iBoWG4J1cJcXGacmvyJigBdSz0dN6PfTw0QILs%2BseAEcOhMa1Kt6AMGB9KvRzGXZE3RUStgG4GJr0JTgOOOt81rEWkEG6chC8OA3oSdWTZ5bRUIdLjTr%2FH0kJuznMzcrgxtQbcx5ZkuXLfPg%2Boc3cTfWQ38L0WWCvc%2BGEw69%2BG6u%2FMGDGK6Je6q3%2BpOyUGt5os68ARaWHLQLz7bgxV58RrEeaAIoQhF83XjN5h8tsTi%2Fo5%2BOslAiMBxAD04nhdnSXXnhVEtCoMJoYGiL9EEimNHvkr7FH%2FvaXAjAGGH5mIoBwPI4Vizs5D4QCbBY0qNJITxj5rmGu6BX6VlO%2FCWLYb588nLajWUmNIffIkUhMOBW9NT%2ByAGb4mBe%2F%2Bpxx4zgEHlENGDwx8oryi17H%2FZG4gfHtlp9rLdlfmv6LSieV7wzpuIuVwXEfog3Tl7K%2Bwlwht28WaVBDDIlfzBBjqmAXV8vljULiWH77CYCOPwGpqMj%2F6Fy43kvX1v5JoUvmJce0YrZ9imG5KMxEmbEiYw9OdvD2UOdkdDk10JL%2FK%2BvxjQPvrnGWEm6ZyzzqW%2FeBkA5GcI7PGqFqnu9vPFsb48GBGfJyEuSNqc%2BMjyVDGn5Qa%2FxSph1mEJRr%2Bn4%2BpDixHIcR2N44v1%2BCFxr5n96pU5zpZDJlZbs%2BiC4SfJICCn8o5wym0nBhI%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2003%20Jun%202025%2023%3A14%3A42%20GMT&X-Amz-Signature=002658dffb502b866e8f56470c81231c33017864cea23c9547143016c2393c3f"
$ws.Range("L2").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1748891365835/MSujKlorgb-camera-video-segment-1748891368335.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSORYA3UK6%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T151442Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjED8aCXVzLXdlc3QtMiJGMEQCIA0KiZttqleQAYFkoPrqGfaQ%2BEfVrMoOu8tcRRo3a5ubAiAR1SS6ywof5CKQTSha4d49IXbz7Aet18sLqGw7V0KBHCrdAwgYEAQaDDc4MTIwNDk0MjI0NCIMKvKB7UOIkRyDvF9OKroDViZbjkBBAv%2FiogeHYKxF%2BDalm8ruSrGzUNAdMd1bL85jwGiCmx%2BvYVaoGRUS6DRdywOkDnUGtHrAWIFtDmLdQgAelEr%2FMOD8N0Z3%2FcKX5oi3cIHVXVh87O5WXRp2Y5kQwrzy4iBoWG4J1cJcXGacmvyJigBdSz0dN6PfTw0QILs%2BseAEcOhMa1Kt6AMGB9KvRzGXZE3RUStgG4GJr0JTgOOOt81rEWkEG6chC8OA3oSdWTZ5bRUIdLjTr%2FH0kJuznMzcrgxtQbcx5ZkuXLfPg%2Boc3cTfWQ38L0WWCvc%2BGEw69%2BG6u%2FMGDGK6Je6q3%2BpOyUGt5os68ARaWHLQLz7bgxV58RrEeaAIoQhF83XjN5h8tsTi%2Fo5%2BOslAiMBxAD04nhdnSXXnhVEtCoMJoYGiL9EEimNHvkr7FH%2FvaXAjAGGH5mIoBwPI4Vizs5D4QCbBY0qNJITxj5rmGu6BX6VlO%2FCWLYb588nLajWUmNIffIkUhMOBW9NT%2ByAGb4mBe%2F%2Bpxx4zgEHlENGDwx8oryi17H%2FZG4gfHtlp9rLdlfmv6LSieV7wzpuIuVwXEfog3Tl7K%2Bwlwht28WaVBDDIlfzBBjqmAXV8vljULiWH77CYCOPwGpqMj%2F6Fy43kvX1v5JoUvmJce0YrZ9imG5KMxEmbEiYw9OdvD2UOdkdDk10JL%2FK%2BvxjQPvrnGWEm6ZyzzqW%2FeBkA5GcI7PGqFqnu9vPFsb48GBGfJyEuSNqc%2BMjyVDGn5Qa%2FxSph1mEJRr%2Bn4%2BpDixHIcR2N44v1%2BCFxr5n96pU5zpZDJlZbs%2BiC4SfJICCn8o5wym0nBhI%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2003%20Jun%202025%2023%3A14%3A42%20GMT&X-Amz-Signature=c52cd41dd10af6bafd0286d30c84f0175b2f60e57bc0a1ba15bdf150e8df2912"

# --- Row 3 ---
$ws.Range("A3").Value = "'281474991205821-1748885146509"
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").Value = "Mobile Usage"
$ws.Range("C3").Value = "'2025-06-02T11:25:46.509"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "'281474991205821"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'148"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "'51834015"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Value = "LUIS IBARRA"
$ws.Range("H3").Value = 20.59709316
$ws.Range("I3").Value = -103.438469289
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205821/1748885144009/VaUCDIqpi8-camera-video-segment-driver-1748885146509.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSORYA3UK6%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T151442Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjED8aCXVzLXdlc3QtMiJGMEQCIA0KiZttqleQAYFkoPrqGfaQ%2BEfVrMoOu8tcRRo3a5ubAiAR1SS6ywof5CKQTSha4d49IXbz7Aet18sLqGw7V0KBHCrdAwgYEAQaDDc4MTIwNDk0MjI0NCIMKvKB7UOIkRyDvF9OKroDViZbjkBBAv%2FiogeHYKxF%2BDalm8ruSrGzUNAdMd1bL85jwGiCmx%2BvYVaoGRUS6DRdywOkDnUGtHrAWIFtDmLdQgAelEr%2FMOD8N0Z3%2FcKX5oi3cIHVXVh87O5WXRp2Y5kQwrzy4iBoWG4J1cJcXGacmvyJigBdSz0dN6PfTw0QILs%2BseAEcOhMa1Kt6AMGB9KvRzGXZE3RUStgG4GJr0JTgOOOt81rEWkEG6chC8OA3oSdWTZ5bRUIdLjTr%2FH0kJuznMzcrgxtQbcx5ZkuXLfPg%2Boc3cTfWQ38L0WWCvc%2BGEw69%2BG6u%2FMGDGK6Je6q3%2BpOyUGt5os68ARaWHLQLz7bgxV58RrEeaAIoQhF83XjN5h8tsTi%2Fo5%2BOslAiMBxAD04nhdnSXXnhVEtCoMJoYGiL9EEimNHvkr7FH%2FvaXAjAGGH5mIoBwPI4Vizs5D4QCbBY0qNJITxj5rmGu6BX6VlO%2FCWLYb588nLajWUmNIffIkUhMOBW9NT%2ByAGb4mBe%2F%2Bpxx4zgEHlENGDwx8oryi17H%2FZG4gfHtlp9rLdlfmv6LSieV7wzpuIuVwXEfog3Tl7K%2Bwlwht28WaVBDDIlfzBBjqmAXV8vljULiWH77CYCOPwGpqMj%2F6Fy43kvX1v5JoUvmJce0YrZ9imG5KMxEmbEiYw9OdvD2UOdkdDk10JL%2FK%2BvxjQPvrnGWEm6ZyzzqW%2FeBkA5GcI7PGqFqnu9vPFsb48GBGfJyEuSNqc%2BMjyVDGn5Qa%2FxSph1mEJRr%2Bn4%2BpDixHIcR2N44v1%2BCFxr5n96pU5zpZDJlZbs%2BiC4SfJICCn8o5wym0nBhI%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2003%20Jun%202025%2023%3A14%3A42%20GMT&X-Amz-Signature=f0e5f5c86abe261610f6d144de021e6b1ea6b45648cd56f21027c73935a71300"
$ws.Range("L3").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205821/1748885144009/8zw63V5M8l-camera-video-segment-1748885146509.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSORYA3UK6%2F20250603%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250603T151442Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjED8aCXVzLXdlc3QtMiJGMEQCIA0KiZttqleQAYFkoPrqGfaQ%2BEfVrMoOu8tcRRo3a5ubAiAR1SS6ywof5CKQTSha4d49IXbz7Aet18sLqGw7V0KBHCrdAwgYEAQaDDc4MTIwNDk0MjI0NCIMKvKB7UOIkRyDvF9OKroDViZbjkBBAv%2FiogeHYKxF%2BDalm8ruSrGzUNAdMd1bL85jwGiCmx%2BvYVaoGRUS6DRdywOkDnUGtHrAWIFtDmLdQgAelEr%2FMOD8N0Z3%2FcKX5oi3cIHVXVh87O5WXRp2Y5kQwrzy4iBoWG4J1cJcXGacmvyJigBdSz0dN6PfTw0QILs%2BseAEcOhMa1Kt6AMGB9KvRzGXZE3RUStgG4GJr0JTgOOOt81rEWkEG6chC8OA3oSdWTZ5bRUIdLjTr%2FH0kJuznMzcrgxtQbcx5ZkuXLfPg%2Boc3cTfWQ38L0WWCvc%2BGEw69%2BG6u%2FMGDGK6Je6q3%2BpOyUGt5os68ARaWHLQLz7bgxV58RrEeaAIoQhF83XjN5h8tsTi%2Fo5%2BOslAiMBxAD04nhdnSXXnhVEtCoMJoYGiL9EEimNHvkr7FH%2FvaXAjAGGH5mIoBwPI4Vizs5D4QCbBY0qNJITxj5rmGu6BX6VlO%2FCWLYb588nLajWUmNIffIkUhMOBW9NT%2ByAGb4mBe%2F%2Bpxx4zgEHlENGDwx8oryi17H%2FZG4gfHtlp9rLdlfmv6LSieV7wzpuIuVwXEfog3Tl7K%2Bwlwht28WaVBDDIlfzBBjqmAXV8vljULiWH77CYCOPwGpqMj%2F6Fy43kvX1v5JoUvmJce0YrZ9imG5KMxEmbEiYw9OdvD2UOdkdDk10JL%2FK%2BvxjQPvrnGWEm6ZyzzqW%2FeBkA5GcI7PGqFqnu9vPFsb48GBGfJyEuSNqc%2BMjyVDGn5Qa%2FxSph1mEJRr%2Bn4%2BpDixHIcR2N44v1%2BCFxr5n96pU5zpZDJlZbs%2BiC4SfJICCn8o5wym0nBhI%3D&X-Amz-SignedHeaders=host&response-expires=Tue%2C%2003%20Jun%202025%2023%3A14%3A42%20GMT&X-Amz-Signature=21f3ced21edc4a3b7d9857421950578b6cc73090a7f3a9e0cdb55724252aff32"
